# Update the "Förändrad" (Changed) date column (C) for every data row,
# and append the "Beteckning" (column A) value as the friendly-name
# second argument of every HYPERLINK() formula on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1

# Data rows start at row 2 (row 1 is the header, row 0 is an empty spacer row).
for ($r = 2; $r -le $lastRow; $r++) {

    # Column C = 3 ("Förändrad" date) -> new serial date value.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -ne $null) {
        $cCell.Value2 = 45186
    }

    # Column A = 1 ("Beteckning") is used as the HYPERLINK friendly name.
    $name = $ws.Cells.Item($r, 1).Value2

    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like 'HYPERLINK(*' -or $f -like '=HYPERLINK(*') {
                if ($f -notlike '*,*') {
                    $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
